# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de), rows 2 and 3 are "handed back":
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - The "Latest Target File" (E) and "Latest Handback File" (F) columns get
#     filled in (same filenames/links as the handoff columns) now that the
#     handback has happened.
#   - "Latest Handback DateTime" (G) is stamped with the handback time.
# The Overview sheet mirrors the same Status text in its zh-cn/de-de columns.

$wb = $excel.ActiveWorkbook

$mdFileName  = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.md"
$statusText  = "Handed back: in sync with en-US"

$zhXlfName = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.zh-cn.xlf"
$deXlfName = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.de-de.xlf"

$mdLink   = "https://github.com/OpenLocalizationTest/oltest/blob/ddd4d8fd425730564d7785a0c6f3eac63e304486/e2e/a9de59e0-3756-49d9-b0b7-ac1152a258f9.md"
$zhXlfLink = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/21e81d0e3f5abf4a9b3eb1cc3fc61ad7e0bc3d7a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.zh-cn.xlf"
$deXlfLink = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/81c6f8aba134159ad79051cc7c8b8c2436809095/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.de-de.xlf"

# ---- Overview sheet: refresh the Status columns for both locales ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B2").Value = $statusText
$zh.Range("E2").Value = $mdFileName
$zh.Hyperlinks.Add($zh.Range("E2"), $mdLink, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdFileName) | Out-Null
$zh.Range("F2").Value = $zhXlfName
$zh.Hyperlinks.Add($zh.Range("F2"), $zhXlfLink, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $zhXlfName) | Out-Null
$zh.Range("G2").Value = "2016-02-17 10:09:30"

$zh.Range("B3").Value = $statusText
$zh.Range("E3").Value = $mdFileName
$zh.Hyperlinks.Add($zh.Range("E3"), $mdLink, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdFileName) | Out-Null
$zh.Range("F3").Value = $zhXlfName
$zh.Hyperlinks.Add($zh.Range("F3"), $zhXlfLink, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $zhXlfName) | Out-Null
$zh.Range("G3").Value = "2016-02-17 10:09:30"

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = $statusText
$de.Range("E2").Value = $mdFileName
$de.Hyperlinks.Add($de.Range("E2"), $mdLink, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdFileName) | Out-Null
$de.Range("F2").Value = $deXlfName
$de.Hyperlinks.Add($de.Range("F2"), $deXlfLink, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $deXlfName) | Out-Null
$de.Range("G2").Value = "2016-02-17 10:09:50"

$de.Range("B3").Value = $statusText
$de.Range("E3").Value = $mdFileName
$de.Hyperlinks.Add($de.Range("E3"), $mdLink, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdFileName) | Out-Null
$de.Range("F3").Value = $deXlfName
$de.Hyperlinks.Add($de.Range("F3"), $deXlfLink, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $deXlfName) | Out-Null
$de.Range("G3").Value = "2016-02-17 10:09:50"

"Handback report generated."
